# Trade #29 closed at 2026-02-16 21:27:32 - leadlag UP +0.000%
#
# Helper: write a literal text value into a cell without Excel's
# "smart" auto-conversion (percentages / dates / numeric-looking
# strings) turning it into a number.
function Set-TextCell {
    param($ws, $row, $col, [string]$text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet: roll up totals now that trades 8-11 closed and
# trade 29 opened.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 3).Value = 11
Set-TextCell $summary 2 4 "63.6%"
Set-TextCell $summary 2 5 "+2.4757%"
Set-TextCell $summary 2 6 "+0.2251%"

$summary.Cells.Item(3, 3).Value = 23
Set-TextCell $summary 3 4 "26.1%"
Set-TextCell $summary 3 5 "+2.4304%"
Set-TextCell $summary 3 6 "+0.1057%"

# ---------------------------------------------------------------
# leadlag sheet: close out trades 8-11 (rows 7-10), then append
# the newly opened trade #29 as row 25.
# ---------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Trade 8 (row 7)
$leadlag.Cells.Item(7, 7).Value = 69082.57287
Set-TextCell $leadlag 7 8 "CLOSED"
$leadlag.Cells.Item(7, 9).Value = 0.4064
$leadlag.Cells.Item(7, 10).Value = 4.06
Set-TextCell $leadlag 7 13 "time_exit_5min"
$leadlag.Cells.Item(7, 14).Value = 5

# Trade 9 (row 8)
$leadlag.Cells.Item(8, 7).Value = 68962.71686099999
Set-TextCell $leadlag 8 8 "CLOSED"
$leadlag.Cells.Item(8, 9).Value = 0.5567
$leadlag.Cells.Item(8, 10).Value = 5.57
Set-TextCell $leadlag 8 13 "time_exit_5min"
$leadlag.Cells.Item(8, 14).Value = 5

# Trade 10 (row 9)
$leadlag.Cells.Item(9, 7).Value = 69394.909679
Set-TextCell $leadlag 9 8 "CLOSED"
$leadlag.Cells.Item(9, 9).Value = -0.1091
$leadlag.Cells.Item(9, 10).Value = -1.09
Set-TextCell $leadlag 9 13 "time_exit_5min"
$leadlag.Cells.Item(9, 14).Value = 5

# Trade 11 (row 10)
$leadlag.Cells.Item(10, 7).Value = 69569.355167
Set-TextCell $leadlag 10 8 "CLOSED"
$leadlag.Cells.Item(10, 9).Value = 0.2389
$leadlag.Cells.Item(10, 10).Value = 2.39
Set-TextCell $leadlag 10 13 "time_exit_5min"
$leadlag.Cells.Item(10, 14).Value = 5

# New trade 29 (row 25) - freshly opened, no exit yet.
$leadlag.Cells.Item(25, 1).Value = 29
Set-TextCell $leadlag 25 2 "2026-02-16"
Set-TextCell $leadlag 25 3 "21:27:32"
Set-TextCell $leadlag 25 4 "leadlag"
Set-TextCell $leadlag 25 5 "UP"
$leadlag.Cells.Item(25, 6).Value = 68972.875
$leadlag.Cells.Item(25, 7).Value = ""
Set-TextCell $leadlag 25 8 "OPEN"
$leadlag.Cells.Item(25, 9).Value = 0
$leadlag.Cells.Item(25, 10).Value = 0
$leadlag.Cells.Item(25, 11).Value = 0.75
Set-TextCell $leadlag 25 12 "Binance leading with 0.078% move"
$leadlag.Cells.Item(25, 13).Value = ""
$leadlag.Cells.Item(25, 14).Value = 0

# ---------------------------------------------------------------
# All Trades sheet: append the four now-closed leadlag trades
# (8, 9, 10, 11) as rows 9-12.
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade 8 -> row 9
$allTrades.Cells.Item(9, 1).Value = 8
Set-TextCell $allTrades 9 2 "2026-02-16"
Set-TextCell $allTrades 9 3 "21:21:59"
Set-TextCell $allTrades 9 4 "leadlag"
Set-TextCell $allTrades 9 5 "DOWN"
$allTrades.Cells.Item(9, 6).Value = 69364.49000000001
$allTrades.Cells.Item(9, 7).Value = 69082.57287
Set-TextCell $allTrades 9 8 "CLOSED"
$allTrades.Cells.Item(9, 9).Value = 0.4064
$allTrades.Cells.Item(9, 10).Value = 4.06
$allTrades.Cells.Item(9, 11).Value = 0.7448
Set-TextCell $allTrades 9 12 "Binance leading with -0.074% move"
Set-TextCell $allTrades 9 13 "time_exit_5min"
$allTrades.Cells.Item(9, 14).Value = 5

# Trade 9 -> row 10
$allTrades.Cells.Item(10, 1).Value = 9
Set-TextCell $allTrades 10 2 "2026-02-16"
Set-TextCell $allTrades 10 3 "21:22:06"
Set-TextCell $allTrades 10 4 "leadlag"
Set-TextCell $allTrades 10 5 "DOWN"
$allTrades.Cells.Item(10, 6).Value = 69348.815
$allTrades.Cells.Item(10, 7).Value = 68962.71686099999
Set-TextCell $allTrades 10 8 "CLOSED"
$allTrades.Cells.Item(10, 9).Value = 0.5567
$allTrades.Cells.Item(10, 10).Value = 5.57
$allTrades.Cells.Item(10, 11).Value = 0.7199
Set-TextCell $allTrades 10 12 "Binance leading with -0.072% move"
Set-TextCell $allTrades 10 13 "time_exit_5min"
$allTrades.Cells.Item(10, 14).Value = 5

# Trade 10 -> row 11
$allTrades.Cells.Item(11, 1).Value = 10
Set-TextCell $allTrades 11 2 "2026-02-16"
Set-TextCell $allTrades 11 3 "21:22:12"
Set-TextCell $allTrades 11 4 "leadlag"
Set-TextCell $allTrades 11 5 "DOWN"
$allTrades.Cells.Item(11, 6).Value = 69319.3
$allTrades.Cells.Item(11, 7).Value = 69394.909679
Set-TextCell $allTrades 11 8 "CLOSED"
$allTrades.Cells.Item(11, 9).Value = -0.1091
$allTrades.Cells.Item(11, 10).Value = -1.09
$allTrades.Cells.Item(11, 11).Value = 0.7119
Set-TextCell $allTrades 11 12 "Binance leading with -0.071% move"
Set-TextCell $allTrades 11 13 "time_exit_5min"
$allTrades.Cells.Item(11, 14).Value = 5

# Trade 11 -> row 12
$allTrades.Cells.Item(12, 1).Value = 11
Set-TextCell $allTrades 12 2 "2026-02-16"
Set-TextCell $allTrades 12 3 "21:22:29"
Set-TextCell $allTrades 12 4 "leadlag"
Set-TextCell $allTrades 12 5 "UP"
$allTrades.Cells.Item(12, 6).Value = 69403.545
$allTrades.Cells.Item(12, 7).Value = 69569.355167
Set-TextCell $allTrades 12 8 "CLOSED"
$allTrades.Cells.Item(12, 9).Value = 0.2389
$allTrades.Cells.Item(12, 10).Value = 2.39
$allTrades.Cells.Item(12, 11).Value = 0.75
Set-TextCell $allTrades 12 12 "Binance leading with 0.134% move"
Set-TextCell $allTrades 12 13 "time_exit_5min"
$allTrades.Cells.Item(12, 14).Value = 5

# ---------------------------------------------------------------
# Comparison sheet: refresh leadlag aggregate stats (row 2).
# ---------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Cells.Item(2, 2).Value = 23
Set-TextCell $comparison 2 3 "26.1%"
Set-TextCell $comparison 2 4 "5.17"
Set-TextCell $comparison 2 5 "+0.5022%"
Set-TextCell $comparison 2 6 "-0.1943%"
Set-TextCell $comparison 2 7 "2.58"
